# Revise and add detail to FSM diagram
# Fix several "SOURCE"/"DEST" state-name labels that were missing a space
# after the leading numeral's period, and normalize "31. ext.ITI" to
# "31. ext_ITI" to match the trigger naming convention used elsewhere.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 10.BeepMinSampling -> 10. BeepMinSampling
$ws.Range("B9").Value  = "10. BeepMinSampling"
$ws.Range("A18").Value = "10. BeepMinSampling"

# 25.timeOut_EarlyWithdraw -> 25. timeOut_EarlyWithdraw
$ws.Range("B10").Value = "25. timeOut_EarlyWithdraw"
$ws.Range("A11").Value = "25. timeOut_EarlyWithdraw"
$ws.Range("B13").Value = "25. timeOut_EarlyWithdraw"

# 30.ITI -> 30. ITI
$ws.Range("B12").Value = "30. ITI"
$ws.Range("B14").Value = "30. ITI"
$ws.Range("B17").Value = "30. ITI"
$ws.Range("A60").Value = "30. ITI"

# 31. ext.ITI -> 31. ext_ITI
$ws.Range("B58").Value = "31. ext_ITI"
$ws.Range("B59").Value = "31. ext_ITI"

# Move the active selection to B6 (matches the saved view state)
$ws.Range("B6").Select()
